$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.897137
$ws.Range("H2").Value = 2.691411
$ws.Range("I2").Value = 0.006112709656408342
$ws.Range("J2").Value = 0.006112709656408343
$ws.Range("M2").Value = 0.366183
$ws.Range("N2").Value = 1.098549
$ws.Range("O2").Value = 0.0639836884691917
$ws.Range("P2").Value = 0.0639836884691917
$ws.Range("Q2").Value = 0.328516318071
$ws.Range("R2").Value = 2.956646862639
$ws.Range("S2").Value = 0.0003911137103582512
$ws.Range("T2").Value = 0.0003911137103582512

# Row 3
$ws.Range("G3").Value = 0.897137
$ws.Range("H3").Value = 2.691411
$ws.Range("I3").Value = 0.006112709656408342
$ws.Range("J3").Value = 0.006112709656408343
$ws.Range("O3").Value = 0.2777364052521014
$ws.Range("P3").Value = 0.2777364052521014
$ws.Range("Q3").Value = 1.426003149093667
$ws.Range("R3").Value = 12.834028341843
$ws.Range("S3").Value = 0.001697722006320661
$ws.Range("T3").Value = 0.001697722006320661

# Row 4
$ws.Range("G4").Value = 0.897137
$ws.Range("H4").Value = 2.691411
$ws.Range("I4").Value = 0.006112709656408342
$ws.Range("J4").Value = 0.006112709656408343
$ws.Range("O4").Value = 0.6582799062787069
$ws.Range("P4").Value = 0.6582799062787069
$ws.Range("Q4").Value = 3.379856589151333
$ws.Range("R4").Value = 30.418709302362
$ws.Range("S4").Value = 0.00402387393972943
$ws.Range("T4").Value = 0.00402387393972943

# Row 5
$ws.Range("I5").Value = 0.9127347171890602
$ws.Range("J5").Value = 0.9127347171890601
$ws.Range("M5").Value = 0.366183
$ws.Range("N5").Value = 1.098549
$ws.Range("O5").Value = 0.0639836884691917
$ws.Range("P5").Value = 0.0639836884691917
$ws.Range("Q5").Value = 49.05324569966701
$ws.Range("R5").Value = 441.4792112970031
$ws.Range("S5").Value = 0.05840013379964062
$ws.Range("T5").Value = 0.0584001337996406

# Row 6
$ws.Range("I6").Value = 0.9127347171890602
$ws.Range("J6").Value = 0.9127347171890601
$ws.Range("O6").Value = 0.2777364052521014
$ws.Range("P6").Value = 0.2777364052521014
$ws.Range("S6").Value = 0.253499659300883
$ws.Range("T6").Value = 0.253499659300883

# Row 7
$ws.Range("I7").Value = 0.9127347171890602
$ws.Range("J7").Value = 0.9127347171890601
$ws.Range("O7").Value = 0.6582799062787069
$ws.Range("P7").Value = 0.6582799062787069
$ws.Range("S7").Value = 0.6008349240885366
$ws.Range("T7").Value = 0.6008349240885366

# Row 8
$ws.Range("I8").Value = 0.08115257315453157
$ws.Range("J8").Value = 0.08115257315453157
$ws.Range("M8").Value = 0.366183
$ws.Range("N8").Value = 1.098549
$ws.Range("O8").Value = 0.0639836884691917
$ws.Range("P8").Value = 0.0639836884691917
$ws.Range("Q8").Value = 4.361395523958
$ws.Range("R8").Value = 39.25255971562201
$ws.Range("S8").Value = 0.005192440959192838
$ws.Range("T8").Value = 0.005192440959192838

# Row 9
$ws.Range("I9").Value = 0.08115257315453157
$ws.Range("J9").Value = 0.08115257315453157
$ws.Range("O9").Value = 0.2777364052521014
$ws.Range("P9").Value = 0.2777364052521014
$ws.Range("S9").Value = 0.02253902394489779
$ws.Range("T9").Value = 0.02253902394489779

# Row 10
$ws.Range("I10").Value = 0.08115257315453157
$ws.Range("J10").Value = 0.08115257315453157
$ws.Range("O10").Value = 0.6582799062787069
$ws.Range("P10").Value = 0.6582799062787069
$ws.Range("S10").Value = 0.05342110825044095
$ws.Range("T10").Value = 0.05342110825044095

